$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 and 17: coin name/link swapped, plus new price & volume values.
# (D16/D17 contain two dots, e.g. "4.139.45" / "98.289.94" - Excel keeps these as text automatically.)
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "4.139.45"
$ws.Range("E16").Value = "  +5.57%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "98.289.94"
$ws.Range("E17").Value = "  +0.88%  "

# Remaining rows: update Price (D) and/or Volume (E) columns.
# "Text" flag marks D values that look like a plain decimal number (a single dot) -
# Excel would silently convert those to a numeric value, so we force the cell to stay
# text (matching the source file's inlineStr cells) and then restore the default
# "Normal" style so no stray style index is left on the cell.
$updates = @(
    @{Row=2;  D="98.787.88";  Text=$false; E="  +1.22%  "},
    @{Row=3;  D="3.483.54";   Text=$false; E="  +5.58%  "},
    @{Row=4;  D=$null;        Text=$false; E="  -0.10%  "},
    @{Row=5;  D="262.58";     Text=$true;  E="  +3.01%  "},
    @{Row=6;  D="679.25";     Text=$true;  E="  +9.54%  "},
    @{Row=7;  D="1.56";       Text=$true;  E="  +9.08%  "},
    @{Row=8;  D=$null;        Text=$false; E="  +16.64%  "},
    @{Row=9;  D=$null;        Text=$false; E="  +22.72%  "},
    @{Row=10; D=$null;        Text=$false; E="  -0.18%  "},
    @{Row=11; D="3.480.76";   Text=$false; E="  +5.56%  "},
    @{Row=12; D=$null;        Text=$false; E="  +11.81%  "},
    @{Row=13; D="42.80";      Text=$true;  E="  +10.93%  "},
    @{Row=14; D=$null;        Text=$false; E="  +12.01%  "},
    @{Row=15; D="6.29";       Text=$true;  E="  +15.41%  "},
    @{Row=18; D="8.17";       Text=$true;  E="  +32.13%  "},
    @{Row=19; D="3.477.34";   Text=$false; E="  +5.46%  "},
    @{Row=20; D="17.63";      Text=$true;  E="  +16.70%  "},
    @{Row=21; D=$null;        Text=$false; E="  +4.10%  "},
    @{Row=22; D="536.94";     Text=$true;  E="  +12.75%  "},
    @{Row=23; D="10.78";      Text=$true;  E="  +15.02%  "},
    @{Row=24; D="0.0000221";  Text=$true;  E="  +8.91%  "},
    @{Row=25; D="0.452";      Text=$true;  E="  +52.04%  "},
    @{Row=26; D="6.29";       Text=$true;  E="  +13.03%  "},
    @{Row=27; D="102.74";     Text=$true;  E="  +16.88%  "},
    @{Row=28; D=$null;        Text=$false; E="  +10.43%  "},
    @{Row=29; D=$null;        Text=$false; E="  +15.29%  "},
    @{Row=30; D=$null;        Text=$false; E="  +8.18%  "},
    @{Row=31; D="11.43";      Text=$true;  E="  +17.04%  "},
    @{Row=32; D=$null;        Text=$false; E="  +0.18%  "},
    @{Row=33; D=$null;        Text=$false; E="  +29.85%  "},
    @{Row=34; D="30.96";      Text=$true;  E="  +12.84%  "},
    @{Row=35; D="0.990";      Text=$true;  E="  -1.08%  "},
    @{Row=36; D=$null;        Text=$false; E="  +16.20%  "},
    @{Row=37; D="8.01";       Text=$true;  E="  +12.33%  "},
    @{Row=38; D=$null;        Text=$false; E="  +10.01%  "},
    @{Row=39; D="542.99";     Text=$true;  E="  +11.25%  "},
    @{Row=40; D="1.43";       Text=$true;  E="  +16.07%  "},
    @{Row=41; D=$null;        Text=$false; E="  -0.06%  "},
    @{Row=42; D="0.872";      Text=$true;  E="  +9.98%  "},
    @{Row=43; D="0.0440";     Text=$true;  E="  +36.18%  "},
    @{Row=44; D="3.52";       Text=$true;  E="  +12.54%  "},
    @{Row=45; D=$null;        Text=$false; E="  +2.52%  "},
    @{Row=46; D="8.25";       Text=$true;  E="  +15.79%  "},
    @{Row=47; D=$null;        Text=$false; E="  +0.01%  "},
    @{Row=48; D="2.14";       Text=$true;  E="  +13.20%  "},
    @{Row=49; D=$null;        Text=$false; E="  +19.49%  "},
    @{Row=50; D=$null;        Text=$false; E="  +15.15%  "},
    @{Row=51; D="51.43";      Text=$true;  E="  +13.36%  "}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        if ($u.Text) {
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}
